# Move the "_GoBack" bookmark from the first of the trailing empty
# paragraphs (right after the diagram paragraph) to the last trailing
# empty paragraph (immediately before the section break / sectPr).

$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark (its paragraph stays, just loses
# the bookmark).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the true last paragraph of the document body (right before the
# section properties). The Paragraphs collection can report spurious
# trailing entries, so find it via the end of the body content instead.
$endPos = $d.Content.End
$lastPara = $d.Range($endPos - 1, $endPos).Paragraphs(1)

$d.Bookmarks.Add("_GoBack", $lastPara.Range)
